# Generate Report for Handoff
# Replace the old handoff UUID / hash-based file names (and the two
# handback timestamps that were regenerated alongside them) with the
# values produced by the new handoff run. Hyperlink *targets* are left
# untouched - only the visible text (cell value + hyperlink display)
# changes, matching the new file names.

$wb = $excel.ActiveWorkbook

$newUuid = "a3dad79a-86b7-40b3-9d55-ec537e1f6203"
$newHash = "d34e5a3279b977f57b25ecac3585a363b692ea23"

$newMdName    = "$newUuid.md"
$newZhXlfName = "$newUuid.$newHash.zh-cn.xlf"
$newDeXlfName = "$newUuid.$newHash.de-de.xlf"

$newZhTimestamp = "2016-03-09 20:58:40"
$newDeTimestamp = "2016-03-09 20:58:45"

# Original hyperlink targets (unchanged by this edit - only the display
# text moves to the new file names).
$mdAddress = "https://github.com/OpenLocalizationTest/oltest/blob/b2700077e475a19ff3d01a8acfadc66f9312644e/e2e/c79a4f70-109c-4362-8d17-02075c1c3b72.md"
$configAddress = "https://github.com/OpenLocalizationTest/oltest/blob/b2700077e475a19ff3d01a8acfadc66f9312644e/.localization-config"
$zhXlfAddress = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/6e977b4f046b801500c62e537bfbc970e153fe0c/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/c79a4f70-109c-4362-8d17-02075c1c3b72.f2d42c1d82ac5c07eb7d46667cf6fc978e33eba6.zh-cn.xlf"
$deXlfAddress = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/645d3c8b61988f8eda0f271d8e2c9057993a541a/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/c79a4f70-109c-4362-8d17-02075c1c3b72.f2d42c1d82ac5c07eb7d46667cf6fc978e33eba6.de-de.xlf"

# ---------------------------------------------------------------------
# Sheet "Overview": only the A2 hyperlink (source .md file) needs its
# display text updated; its target keeps pointing at the same commit.
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("A2"), $mdAddress, "", "", $newMdName)
$wsOverview.Hyperlinks.Add($wsOverview.Range("A3"), $configAddress, "", "", ".localization-config")

# ---------------------------------------------------------------------
# Sheet "zh-cn": A2 (.md) + C2 (.xlf) hyperlinks, plus the D2 handback
# timestamp cell.
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Hyperlinks.Delete()
$wsZh.Hyperlinks.Add($wsZh.Range("A2"), $mdAddress, "", "", $newMdName)
$wsZh.Hyperlinks.Add($wsZh.Range("C2"), $zhXlfAddress, "", "", $newZhXlfName)
$wsZh.Hyperlinks.Add($wsZh.Range("A3"), $configAddress, "", "", ".localization-config")

$wsZh.Range("D2").Value = $newZhTimestamp

# ---------------------------------------------------------------------
# Sheet "de-de": A2 (.md) + C2 (.xlf) hyperlinks, plus the D2 handback
# timestamp cell.
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Hyperlinks.Delete()
$wsDe.Hyperlinks.Add($wsDe.Range("A2"), $mdAddress, "", "", $newMdName)
$wsDe.Hyperlinks.Add($wsDe.Range("C2"), $deXlfAddress, "", "", $newDeXlfName)
$wsDe.Hyperlinks.Add($wsDe.Range("A3"), $configAddress, "", "", ".localization-config")

$wsDe.Range("D2").Value = $newDeTimestamp
